# Update RelTol/AbsTol benchmark values so that all three simulation
# methods (ODE15s, CVODE(SundialsTB), CVODE(IQMTools)) report the same
# RelTol / AbsTol results in the Summary sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("diary")

$ws.Cells.Item(2, 3).Value = 0.0671
$ws.Cells.Item(2, 4).Value = 0.00289

$ws.Cells.Item(3, 3).Value = 0.0261
$ws.Cells.Item(3, 4).Value = 0.000968

$ws.Cells.Item(4, 3).Value = 0.00241
$ws.Cells.Item(4, 4).Value = 0.000169

$ws.Cells.Item(5, 3).Value = 0.0903
$ws.Cells.Item(5, 4).Value = 0.00315

$ws.Cells.Item(6, 3).Value = 0.0334
$ws.Cells.Item(6, 4).Value = 0.00118

$ws.Cells.Item(7, 3).Value = 0.00262
$ws.Cells.Item(7, 4).Value = 0.000226

$ws.Cells.Item(8, 3).Value = 0.0956
$ws.Cells.Item(8, 4).Value = 0.0179

$ws.Cells.Item(9, 3).Value = 0.0689
$ws.Cells.Item(9, 4).Value = 0.00188

$ws.Cells.Item(10, 3).Value = 0.00388
$ws.Cells.Item(10, 4).Value = 0.000135

$ws.Cells.Item(11, 3).Value = 0.115
$ws.Cells.Item(11, 4).Value = 0.0307

$ws.Cells.Item(12, 3).Value = 0.0838
$ws.Cells.Item(12, 4).Value = 0.00197

$ws.Cells.Item(13, 3).Value = 0.00309
$ws.Cells.Item(13, 4).Value = 0.00016

$ws.Cells.Item(14, 3).Value = 0.0321
$ws.Cells.Item(14, 4).Value = 0.00824

$ws.Cells.Item(15, 3).Value = 0.0199
$ws.Cells.Item(15, 4).Value = 0.000687

$ws.Cells.Item(16, 3).Value = 0.00222
$ws.Cells.Item(16, 4).Value = 0.000118

$ws.Cells.Item(17, 3).Value = 0.0159
$ws.Cells.Item(17, 4).Value = 0.00503

$ws.Cells.Item(18, 3).Value = 0.0104
$ws.Cells.Item(18, 4).Value = 0.000491

$ws.Cells.Item(19, 3).Value = 0.0021
$ws.Cells.Item(19, 4).Value = 0.000121

$wb.Save()
